$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($r in 2..13) {
    $ws.Cells.Item($r, 3).Value = 45212
}
